$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 12.07002266666666
$ws.Range("H2").Value = 36.21006799999999
$ws.Range("I2").Value = 0.7601982364861632
$ws.Range("J2").Value = 0.7601982364861634
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 21.09934133333334
$ws.Range("N2").Value = 63.29802400000001
$ws.Range("O2").Value = 0.2917236204149438
$ws.Range("P2").Value = 0.2917236204149438
$ws.Range("Q2").Value = 254.6695281450702
$ws.Range("R2").Value = 2292.025753305631
$ws.Range("S2").Value = 0.2217677817807991
$ws.Range("T2").Value = 0.2217677817807992
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 12.07002266666666
$ws.Range("H3").Value = 36.21006799999999
$ws.Range("I3").Value = 0.7601982364861632
$ws.Range("J3").Value = 0.7601982364861634
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 35.81943766666667
$ws.Range("N3").Value = 107.458313
$ws.Range("O3").Value = 0.4952465516465762
$ws.Range("P3").Value = 0.4952465516465762
$ws.Range("Q3").Value = 432.3414245439204
$ws.Range("R3").Value = 3891.072820895283
$ws.Range("S3").Value = 0.3764855551875808
$ws.Range("T3").Value = 0.3764855551875809
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 12.07002266666666
$ws.Range("H4").Value = 36.21006799999999
$ws.Range("I4").Value = 0.7601982364861632
$ws.Range("J4").Value = 0.7601982364861634
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 15.40769666666667
$ws.Range("N4").Value = 46.22309
$ws.Range("O4").Value = 0.2130298279384801
$ws.Range("P4").Value = 0.2130298279384801
$ws.Range("Q4").Value = 185.9712480077911
$ws.Range("R4").Value = 1673.74123207012
$ws.Range("S4").Value = 0.1619448995177833
$ws.Range("T4").Value = 0.1619448995177834
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.308268
$ws.Range("H5").Value = 3.924804
$ws.Range("I5").Value = 0.08239777620284613
$ws.Range("J5").Value = 0.08239777620284613
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 21.09934133333334
$ws.Range("N5").Value = 63.29802400000001
$ws.Range("O5").Value = 0.2917236204149438
$ws.Range("P5").Value = 0.2917236204149438
$ws.Range("Q5").Value = 27.60359308747734
$ws.Range("R5").Value = 248.432337787296
$ws.Range("S5").Value = 0.02403737758803457
$ws.Range("T5").Value = 0.02403737758803458
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.308268
$ws.Range("H6").Value = 3.924804
$ws.Range("I6").Value = 0.08239777620284613
$ws.Range("J6").Value = 0.08239777620284613
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 35.81943766666667
$ws.Range("N6").Value = 107.458313
$ws.Range("O6").Value = 0.4952465516465762
$ws.Range("P6").Value = 0.4952465516465762
$ws.Range("Q6").Value = 46.86142407729466
$ws.Range("R6").Value = 421.752816695652
$ws.Range("S6").Value = 0.04080721452780586
$ws.Range("T6").Value = 0.04080721452780587
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.308268
$ws.Range("H7").Value = 3.924804
$ws.Range("I7").Value = 0.08239777620284613
$ws.Range("J7").Value = 0.08239777620284613
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 15.40769666666667
$ws.Range("N7").Value = 46.22309
$ws.Range("O7").Value = 0.2130298279384801
$ws.Range("P7").Value = 0.2130298279384801
$ws.Range("Q7").Value = 20.15739650270667
$ws.Range("R7").Value = 181.41656852436
$ws.Range("S7").Value = 0.0175531840870057
$ws.Range("T7").Value = 0.0175531840870057
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 2.499176666666667
$ws.Range("H8").Value = 7.49753
$ws.Range("I8").Value = 0.1574039873109905
$ws.Range("J8").Value = 0.1574039873109906
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 21.09934133333334
$ws.Range("N8").Value = 63.29802400000001
$ws.Range("O8").Value = 0.2917236204149438
$ws.Range("P8").Value = 0.2917236204149438
$ws.Range("Q8").Value = 52.73098154230222
$ws.Range("R8").Value = 474.5788338807201
$ws.Range("S8").Value = 0.04591846104611003
$ws.Range("T8").Value = 0.04591846104611005
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 2.499176666666667
$ws.Range("H9").Value = 7.49753
$ws.Range("I9").Value = 0.1574039873109905
$ws.Range("J9").Value = 0.1574039873109906
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 35.81943766666667
$ws.Range("N9").Value = 107.458313
$ws.Range("O9").Value = 0.4952465516465762
$ws.Range("P9").Value = 0.4952465516465762
$ws.Range("Q9").Value = 89.51910282965444
$ws.Range("R9").Value = 805.67192546689
$ws.Range("S9").Value = 0.07795378193118949
$ws.Range("T9").Value = 0.07795378193118951
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 2.499176666666667
$ws.Range("H10").Value = 7.49753
$ws.Range("I10").Value = 0.1574039873109905
$ws.Range("J10").Value = 0.1574039873109906
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 15.40769666666667
$ws.Range("N10").Value = 46.22309
$ws.Range("O10").Value = 0.2130298279384801
$ws.Range("P10").Value = 0.2130298279384801
$ws.Range("Q10").Value = 38.50655599641111
$ws.Range("R10").Value = 346.5590039677
$ws.Range("S10").Value = 0.03353174433369101
$ws.Range("T10").Value = 0.03353174433369103
